{"js": "// Replace the Jinja2 expression `today(format='...')` with\n// `format_date(rechnungsdatum, format='...')` in the \"Hamburg, {{ ... }}\" line.\nconst body = context.document.body;\n\nconst results = body.search(\"{{ today(format=\", { matchCase: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Target text '{{ today(format=' not found in document body.\");\n}\n\nfor (let i = 0; i < results.items.length; i++) {\n  const range = results.items[i];\n  range.insertText(\"{{ format_date(rechnungsdatum, format=\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the Jinja2 expression `today(format='...')` with\n# `format_date(rechnungsdatum, format='...')` in the \"Hamburg, {{ ... }}\" line.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = \"{{ today(format=\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"{{ format_date(rechnungsdatum, format=\"\n$find.Forward = $true\n$find.Wrap = 1\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.Execute([ref]$find.Text, [ref]$find.MatchCase, [ref]$find.MatchWholeWord, $null, $null, $null, [ref]$find.Forward, [ref]$find.Wrap, $null, [ref]$find.Replacement.Text, 2)\n"}
